# Add 2022-Q4 data:
#   - new worksheet "2022-Q4" (fund holding detail), inserted right after "总计"
#   - new summary row for 2022-Q4 at the top of the "总计" sheet's data
#     (existing 2022-Q3 / 2022-Q2 / 2022-Q1 rows shift down by one)

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

function Set-TextCell($sheet, $row, $col, $val) {
    # Force the cell to be stored as text (not auto-coerced to a number),
    # and then drop back to the "Normal" style so no stray numFmt/quotePrefix
    # style index gets attached to the cell.
    $c = $sheet.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Reuse the header + "A" id-column formatting (bold / centered / bordered,
# style index 2) from the existing "2022-Q3" sheet, which has the exact same
# layout.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("B1:H1").Copy($q4.Range("B1:H1"))
$q3.Range("A2:A4").Copy($q4.Range("A2:A6"))

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Row 2
$q4.Cells.Item(2,1).Value = 0
Set-TextCell $q4 2 2 "159792"
Set-TextCell $q4 2 3 "富国中证港股通互联网ETF"
Set-TextCell $q4 2 4 "38.40"
Set-TextCell $q4 2 5 "99.63"
Set-TextCell $q4 2 6 "2.75"
Set-TextCell $q4 2 7 "1.0560"
$q4.Cells.Item(2,8).Value = 10

# Row 3
$q4.Cells.Item(3,1).Value = 1
Set-TextCell $q4 3 2 "005228"
Set-TextCell $q4 3 3 "汇添富港股通专注成长混合"
Set-TextCell $q4 3 4 "7.51"
Set-TextCell $q4 3 5 "85.64"
Set-TextCell $q4 3 6 "3.81"
Set-TextCell $q4 3 7 "0.2861"
$q4.Cells.Item(3,8).Value = 9

# Row 4
$q4.Cells.Item(4,1).Value = 2
Set-TextCell $q4 4 2 "013123"
Set-TextCell $q4 4 3 "汇添富精选核心优势一年持有混合A"
Set-TextCell $q4 4 4 "5.58"
Set-TextCell $q4 4 5 "83.43"
Set-TextCell $q4 4 6 "3.20"
Set-TextCell $q4 4 7 "0.1786"
$q4.Cells.Item(4,8).Value = 9

# Row 5
$q4.Cells.Item(5,1).Value = 3
Set-TextCell $q4 5 2 "513770"
Set-TextCell $q4 5 3 "华宝中证港股通互联网ETF"
Set-TextCell $q4 5 4 "5.89"
Set-TextCell $q4 5 5 "98.59"
Set-TextCell $q4 5 6 "2.73"
Set-TextCell $q4 5 7 "0.1608"
$q4.Cells.Item(5,8).Value = 10

# Row 6
$q4.Cells.Item(6,1).Value = 4
Set-TextCell $q4 6 2 "013124"
Set-TextCell $q4 6 3 "汇添富精选核心优势一年持有混合C"
Set-TextCell $q4 6 4 "0.27"
Set-TextCell $q4 6 5 "83.43"
Set-TextCell $q4 6 6 "3.20"
Set-TextCell $q4 6 7 "0.0086"
$q4.Cells.Item(6,8).Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2022-Q3 / 2022-Q2 / 2022-Q1
#    rows down by one and insert the new 2022-Q4 row at the top.
# ---------------------------------------------------------------------------

# Extend formatting for the new bottom row (row 5) by copying row 4's
# formats+values down first; the values are overwritten right after.
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))

# Row 5: 2022-Q1 (previously row 4)
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(5,2).Value = "2022-Q1"
$totalSheet.Cells.Item(5,3).Value = 1
$totalSheet.Cells.Item(5,4).Value = 0

# Row 4: 2022-Q2 (previously row 3)
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2022-Q2"
$totalSheet.Cells.Item(4,3).Value = 6
$totalSheet.Cells.Item(4,4).Value = 0.52

# Row 3: 2022-Q3 (previously row 2)
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q3"
$totalSheet.Cells.Item(3,3).Value = 3
$totalSheet.Cells.Item(3,4).Value = 0.91

# Row 2: 2022-Q4 (new)
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 5
$totalSheet.Cells.Item(2,4).Value = 1.69

# Keep "总计" the active sheet/tab, matching the workbook's unchanged
# bookViews (activeTab stays on the first sheet).
$totalSheet.Activate()
